$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing column F (GFA - Sales),
# shifting the old F:M data to H:O.
$ws.Range("F1:G1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("F1").Value = "M_TotalTax"
$ws.Range("G1").Value = "M_CorpTax"

# Copy the header style (bold, bordered, centered) from an existing header
# cell onto the two new header cells.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New numeric data for rows 2-11 (M_TotalTax, M_CorpTax).
$values = @{
    2  = @(6308727034979.312, 399825921028.5854)
    3  = @(16630145391623.02, 1639742485782.957)
    4  = @(4450994137606.095, 601350231413.5104)
    5  = @(4183547438952.192, 598849276038.3025)
    6  = @(11223287075501.79, 872292028558.4308)
    7  = @(1841737275230.086, 214321200777.9413)
    8  = @(6192585801479.285, 516695167857.3162)
    9  = @(14653861967257.56, 1232540278767.842)
    10 = @(9623160693235.053, 876943418066.7275)
    11 = @(5030701274022.499, 355596860701.1148)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 6).Value = $pair[0]
    $ws.Cells.Item($row, 7).Value = $pair[1]
}
